# Spiel_Laden_Fenster_Beschreibung.xlsx - "Überarbeitete Dialogskizzen und Dialogbeschreibungen"
#
# Adds two new description rows (Ausgewählte Datei / Benötigte Mitspieler)
# to the dialog description table and extends the "Datei Auswählen" action
# text to mention that the "Benötigte Mitspieler" label gets populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 7: "Ausgewählte Datei" label description -----------------
$ws.Range("A7").Value2 = "Ausgewählte Datei"
$ws.Range("B7").Value2 = "Label"
$ws.Range("C7").Value2 = "Name der ausgewählten Datei"
$ws.Rows.Item(7).RowHeight = 31.5

# --- New row 8: "Benötigte Mitspieler" label description ---------------
$ws.Range("A8").Value2 = "Benötigte Mitspieler"
$ws.Range("B8").Value2 = "Label"
$ws.Range("C8").Value2 = "Namen der Benötigten Mitspieler"
$ws.Rows.Item(8).RowHeight = 47.25

# --- Extend the "Datei Auswählen" action description (row 5) -----------
$ws.Range("D5").Value2 = "Es wird ein Explorer Fenster geöffnet um eine Speicherdatei zum laden auszuwählen. Nach dem auswählen wird das Label benötigte Mitspieler aus dieser Datei befüllt"
$ws.Rows.Item(5).RowHeight = 126

# --- Update the current selection/view state ----------------------------
[void]$ws.Range("D6").Select()
